$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object 'object[,]' 12,19
$arr[0,0] = 0.98
$arr[0,1] = 0.86
$arr[0,2] = 0.86
$arr[0,3] = 0.79
$arr[0,4] = 0.68
$arr[0,5] = 0.68
$arr[0,6] = 0.5608663317
$arr[0,7] = 0.57
$arr[0,8] = 0.92
$arr[0,9] = 0.5600000000000001
$arr[0,10] = 0.5600000000000001
$arr[0,11] = 0.83
$arr[0,12] = 0.5499221106
$arr[0,13] = 0.3898693467
$arr[0,14] = 0.359718593
$arr[0,15] = 0.4697160804
$arr[0,16] = 0.76
$arr[0,17] = 0.68
$arr[0,18] = 0.73
$arr[1,0] = 0.98
$arr[1,1] = 0.84
$arr[1,2] = 0.86
$arr[1,3] = 0.83
$arr[1,4] = 0.77
$arr[1,5] = 0.76
$arr[1,6] = 0.5806613065
$arr[1,7] = 0.59
$arr[1,8] = 0.93
$arr[1,9] = 0.53
$arr[1,10] = 0.5200829146
$arr[1,11] = 0.9
$arr[1,12] = 0.6000552764
$arr[1,13] = 0.3598462312
$arr[1,14] = 0.3796934673
$arr[1,15] = 0.5
$arr[1,16] = 0.76
$arr[1,17] = 0.62
$arr[1,18] = 0.66
$arr[2,0] = 1
$arr[2,1] = 0.8100000000000001
$arr[2,2] = 0.8
$arr[2,3] = 0.85
$arr[2,4] = 0.78
$arr[2,5] = 0.760801005
$arr[2,6] = 0.6005386935
$arr[2,7] = 0.64
$arr[2,8] = 0.93
$arr[2,9] = 0.53
$arr[2,10] = 0.53
$arr[2,11] = 0.85
$arr[2,12] = 0.5599798995
$arr[2,13] = 0.3498894472
$arr[2,14] = 0.3698341709
$arr[2,15] = 0.41
$arr[2,16] = 0.7
$arr[2,17] = 0.67
$arr[2,18] = 0.73
$arr[3,0] = 0.99
$arr[3,1] = 0.79
$arr[3,2] = 0.8
$arr[3,3] = 0.77
$arr[3,4] = 0.73
$arr[3,5] = 0.71
$arr[3,6] = 0.6008030151
$arr[3,7] = 0.62
$arr[3,8] = 0.92
$arr[3,9] = 0.53
$arr[3,10] = 0.53
$arr[3,11] = 0.8100000000000001
$arr[3,12] = 0.4699020101
$arr[3,13] = 0.3198115578
$arr[3,14] = 0.3898040201
$arr[3,15] = 0.38
$arr[3,16] = 0.65
$arr[3,17] = 0.59
$arr[3,18] = 0.65
$arr[4,0] = 0.99
$arr[4,1] = 0.83
$arr[4,2] = 0.83
$arr[4,3] = 0.78
$arr[4,4] = 0.71
$arr[4,5] = 0.71
$arr[4,6] = 0.58
$arr[4,7] = 0.5700663317
$arr[4,8] = 0.96
$arr[4,9] = 0.52
$arr[4,10] = 0.5000562814
$arr[4,11] = 0.82
$arr[4,12] = 0.5498894472
$arr[4,13] = 0.3799170854
$arr[4,14] = 0.3096879397
$arr[4,15] = 0.43
$arr[4,16] = 0.8
$arr[4,17] = 0.68
$arr[4,18] = 0.6910050250999999
$arr[5,0] = 1
$arr[5,1] = 0.84
$arr[5,2] = 0.84
$arr[5,3] = 0.84
$arr[5,4] = 0.75
$arr[5,5] = 0.7408919598
$arr[5,6] = 0.650741206
$arr[5,7] = 0.6601854271000001
$arr[5,8] = 0.93
$arr[5,9] = 0.55
$arr[5,10] = 0.5301482412
$arr[5,11] = 0.84
$arr[5,12] = 0.4899854271
$arr[5,13] = 0.2898326633
$arr[5,14] = 0.3596939698
$arr[5,15] = 0.3895673367
$arr[5,16] = 0.65
$arr[5,17] = 0.67
$arr[5,18] = 0.7009487437
$arr[6,0] = 1
$arr[6,1] = 0.79
$arr[6,2] = 0.78
$arr[6,3] = 0.75
$arr[6,4] = 0.65
$arr[6,5] = 0.64
$arr[6,6] = 0.5008080402
$arr[6,7] = 0.53
$arr[6,8] = 0.93
$arr[6,9] = 0.44
$arr[6,10] = 0.44
$arr[6,11] = 0.8
$arr[6,12] = 0.5299949749
$arr[6,13] = 0.3900623116
$arr[6,14] = 0.319678392
$arr[6,15] = 0.3395376884
$arr[6,16] = 0.71
$arr[6,17] = 0.6899999999999999
$arr[6,18] = 0.6908994975
$arr[7,0] = 1
$arr[7,1] = 0.83
$arr[7,2] = 0.83
$arr[7,3] = 0.8100000000000001
$arr[7,4] = 0.77
$arr[7,5] = 0.76
$arr[7,6] = 0.6507075377
$arr[7,7] = 0.67
$arr[7,8] = 0.91
$arr[7,9] = 0.54
$arr[7,10] = 0.5201165829
$arr[7,11] = 0.84
$arr[7,12] = 0.5299356784
$arr[7,13] = 0.3197698492
$arr[7,14] = 0.2896190955
$arr[7,15] = 0.3695763819
$arr[7,16] = 0.74
$arr[7,17] = 0.67
$arr[7,18] = 0.7
$arr[8,0] = 0.99
$arr[8,1] = 0.88
$arr[8,2] = 0.89
$arr[8,3] = 0.77
$arr[8,4] = 0.64
$arr[8,5] = 0.63
$arr[8,6] = 0.5409723618
$arr[8,7] = 0.5401984925
$arr[8,8] = 0.91
$arr[8,9] = 0.5
$arr[8,10] = 0.5
$arr[8,11] = 0.8100000000000001
$arr[8,12] = 0.5500844221
$arr[8,13] = 0.3398080402
$arr[8,14] = 0.3696829146
$arr[8,15] = 0.4797361809
$arr[8,16] = 0.72
$arr[8,17] = 0.7
$arr[8,18] = 0.6908994975
$arr[9,0] = 0.99
$arr[9,1] = 0.83
$arr[9,2] = 0.84
$arr[9,3] = 0.77
$arr[9,4] = 0.72
$arr[9,5] = 0.7109266332
$arr[9,6] = 0.5106552764
$arr[9,7] = 0.52
$arr[9,8] = 0.95
$arr[9,9] = 0.49
$arr[9,10] = 0.4699733668
$arr[9,11] = 0.83
$arr[9,12] = 0.449841206
$arr[9,13] = 0.2798030151
$arr[9,14] = 0.2997603015
$arr[9,15] = 0.38
$arr[9,16] = 0.73
$arr[9,17] = 0.68
$arr[9,18] = 0.7109266332
$arr[10,0] = 1
$arr[10,1] = 0.79
$arr[10,2] = 0.79
$arr[10,3] = 0.8
$arr[10,4] = 0.74
$arr[10,5] = 0.74
$arr[10,6] = 0.6307859296
$arr[10,7] = 0.64
$arr[10,8] = 0.95
$arr[10,9] = 0.5600000000000001
$arr[10,10] = 0.5600000000000001
$arr[10,11] = 0.88
$arr[10,12] = 0.479961809
$arr[10,13] = 0.3298738693
$arr[10,14] = 0.369839196
$arr[10,15] = 0.4
$arr[10,16] = 0.78
$arr[10,17] = 0.63
$arr[10,18] = 0.65
$arr[11,0] = 1
$arr[11,1] = 0.82
$arr[11,2] = 0.8
$arr[11,3] = 0.76
$arr[11,4] = 0.68
$arr[11,5] = 0.68
$arr[11,6] = 0.6009
$arr[11,7] = 0.6002869347000001
$arr[11,8] = 0.89
$arr[11,9] = 0.52
$arr[11,10] = 0.5101326633
$arr[11,11] = 0.85
$arr[11,12] = 0.5200241206
$arr[11,13] = 0.349760804
$arr[11,14] = 0.3897477387
$arr[11,15] = 0.4
$arr[11,16] = 0.71
$arr[11,17] = 0.72
$arr[11,18] = 0.7807834171

$ws.Range("A90:S101").Value = $arr

Write-Output "Rows added successfully"